$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Clear the old data block first (values become fully dereferenced),
# then retype every cell row-by-row in the new final order so the
# shared-string table gets rebuilt in first-use order.
$ws.Range("B16:G19").ClearContents()

$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1235045954"
$ws.Range("D16").Value = "CESAR ANDRES DUEÑAS D ELIA"
$ws.Range("E16").Value = "1805"
$ws.Range("F16").Value = 44000
$ws.Range("G16").Value = 1100000

$ws.Range("B17").Value = "PE"
$ws.Range("C17").Value = "927514822121998"
$ws.Range("D17").Value = "ANGEL ALEXANDER PELOCHE TANG"
$ws.Range("E17").Value = "1805"
$ws.Range("F17").Value = 40000
$ws.Range("G17").Value = 1000000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1235045954"
$ws.Range("D18").Value = "CESAR ANDRES DUEÑAS D ELIA"
$ws.Range("E18").Value = "1806"
$ws.Range("F18").Value = 44000
$ws.Range("G18").Value = 1100000

$ws.Range("B19").Value = "PE"
$ws.Range("C19").Value = "927514822121998"
$ws.Range("D19").Value = "ANGEL ALEXANDER PELOCHE TANG"
$ws.Range("E19").Value = "1806"
$ws.Range("F19").Value = 40000
$ws.Range("G19").Value = 1000000
